$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.026.00'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '1.831.50'
$ws.Range('E3').Value = '  -0.62%  '
$__style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = $__style
$ws.Range('E4').Value = '  -0.08%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.15'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  -0.20%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6244'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  -5.75%  '
$__style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = $__style
$ws.Range('E7').Value = '  -0.05%  '
$__style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07576'
$ws.Range('D8').Style = $__style
$ws.Range('E8').Value = '  +1.58%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2917'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  -1.42%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.51'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  -3.30%  '
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').Value = '1.829.14'
$ws.Range('E12').Value = '  -0.78%  '
$__style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.948'
$ws.Range('D13').Style = $__style
$ws.Range('E13').Value = '  -1.42%  '
$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6631'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  -1.62%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001009'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  +15.59%  '
$__style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.61'
$ws.Range('D16').Style = $__style
$__style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.026'
$ws.Range('D17').Style = $__style
$ws.Range('E17').Value = '  -2.45%  '
$ws.Range('D18').Value = '28.988.75'
$ws.Range('E18').Value = '  -0.74%  '
$__style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '226.15'
$ws.Range('D19').Style = $__style
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('E20').Value = '  -1.55%  '
$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9993'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  -0.16%  '
$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.151'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  -0.58%  '
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  -0.06%  '
$__style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.06'
$ws.Range('D24').Style = $__style
$ws.Range('E24').Value = '  -0.57%  '
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.459'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  -2.04%  '
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('E27').Value = '  -0.79%  '
$__style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.491'
$ws.Range('D28').Style = $__style
$ws.Range('E28').Value = '  -1.50%  '
$ws.Range('E29').Value = '  -1.24%  '
$__style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.017'
$ws.Range('D30').Style = $__style
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('E31').Value = '  -1.03%  '
$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05197'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  -3.42%  '
$__style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.842'
$ws.Range('D33').Style = $__style
$ws.Range('E33').Value = '  -0.81%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7361'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  -1.44%  '
$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.138'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('E36').Value = '  +1.76%  '
$ws.Range('D37').Value = '1.244.07'
$ws.Range('E37').Value = '  -4.29%  '
$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.756'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  -0.13%  '
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01781'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  -0.89%  '
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.317'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  -0.71%  '
$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8946'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  -1.03%  '
$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  +0.03%  '
$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.29'
$ws.Range('D43').Style = $__style
$ws.Range('E43').Value = '  -2.32%  '
$ws.Range('D44').Value = '1.978.34'
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('E45').Value = '  +1.82%  '
$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.91'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  -1.72%  '
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4022'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  -0.11%  '
$__style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.826'
$ws.Range('D49').Style = $__style
$ws.Range('E49').Value = '  -0.06%  '
$__style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05750'
$ws.Range('D50').Style = $__style
$ws.Range('E50').Value = '  -1.98%  '
$__style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.633'
$ws.Range('D51').Style = $__style
$ws.Range('E51').Value = '  -6.85%  '
